$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3837.25
$ws.Range("I86").Value = 3949.6667
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 3949.6667
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = -2826.6667
$ws.Range("N86").Value = -5746

$ws.Range("H89").Value = 3837.25
$ws.Range("I89").Value = 3949.6667
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 19748.3335
$ws.Range("L89").Value = 17500
$ws.Range("M89").Value = -14132.3335
$ws.Range("N89").Value = -28732

$ws.Range("H125").Value = 855296
$ws.Range("J125").Value = 832.25
$ws.Range("L125").Value = 7490.25
$ws.Range("N125").Value = -12410.25

$ws.Range("H137").Value = 8535.546
$ws.Range("I137").Value = 14920.2
$ws.Range("K137").Value = 44760.60000000001
$ws.Range("M137").Value = -42210.60000000001

$ws.Range("H138").Value = 8779.634
$ws.Range("J138").Value = 7202.5806
$ws.Range("L138").Value = 21607.7418
$ws.Range("N138").Value = -31887.7418


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32012.75
$ws.Range("I32").Value = 30170.656
$ws.Range("K32").Value = 30170.656
$ws.Range("M32").Value = -29883.656

$ws.Range("H45").Value = 4396.643
$ws.Range("J45").Value = 5617.875
$ws.Range("L45").Value = 5617.875
$ws.Range("N45").Value = -6371.875

$ws.Range("H63").Value = 3624.4119
$ws.Range("I63").Value = 1826.125
$ws.Range("J63").Value = 5222.8887
$ws.Range("K63").Value = 1826.125
$ws.Range("L63").Value = 5222.8887
$ws.Range("M63").Value = -1140.125
$ws.Range("N63").Value = -6594.8887

$ws.Range("H66").Value = 3624.4119
$ws.Range("I66").Value = 1826.125
$ws.Range("J66").Value = 5222.8887
$ws.Range("K66").Value = 9130.625
$ws.Range("L66").Value = 26114.4435
$ws.Range("M66").Value = -5698.625
$ws.Range("N66").Value = -32978.4435

$ws.Range("H74").Value = 3362.9375
$ws.Range("I74").Value = 3686.8572
$ws.Range("K74").Value = 3686.8572
$ws.Range("M74").Value = -2812.8572

$ws.Range("H77").Value = 3362.9375
$ws.Range("I77").Value = 3686.8572
$ws.Range("K77").Value = 18434.286
$ws.Range("M77").Value = -14066.286


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1424
$ws.Range("I86").Value = 1006
$ws.Range("K86").Value = 1006
$ws.Range("M86").Value = 117

$ws.Range("H89").Value = 1424
$ws.Range("I89").Value = 1006
$ws.Range("K89").Value = 5030
$ws.Range("M89").Value = 586


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6490.923
$ws.Range("I31").Value = 4548.2
$ws.Range("K31").Value = 4548.2
$ws.Range("M31").Value = -4253.2

$ws.Range("H34").Value = 6490.923
$ws.Range("I34").Value = 4548.2
$ws.Range("K34").Value = 4548.2
$ws.Range("M34").Value = -4346.2

$ws.Range("H58").Value = 4836.16
$ws.Range("I58").Value = 3954.818
$ws.Range("J58").Value = 5528.643
$ws.Range("K58").Value = 3954.818
$ws.Range("L58").Value = 5528.643
$ws.Range("M58").Value = -3751.818
$ws.Range("N58").Value = -5934.643

$ws.Range("H59").Value = 41497.5
$ws.Range("J59").Value = 41996.668
$ws.Range("L59").Value = 41996.668
$ws.Range("N59").Value = -44286.668

$ws.Range("H86").Value = 3334
$ws.Range("I86").Value = 2248.5
$ws.Range("J86").Value = 3954.2856
$ws.Range("K86").Value = 2248.5
$ws.Range("L86").Value = 3954.2856
$ws.Range("M86").Value = -1125.5
$ws.Range("N86").Value = -6200.2856

$ws.Range("H89").Value = 3334
$ws.Range("I89").Value = 2248.5
$ws.Range("J89").Value = 3954.2856
$ws.Range("K89").Value = 11242.5
$ws.Range("L89").Value = 19771.428
$ws.Range("M89").Value = -5626.5
$ws.Range("N89").Value = -31003.428

$ws.Range("H134").Value = 2330.8538
$ws.Range("I134").Value = 1669.3243
$ws.Range("K134").Value = 5007.9729
$ws.Range("M134").Value = -2472.9729

$ws.Range("H136").Value = 4836.16
$ws.Range("I136").Value = 3954.818
$ws.Range("J136").Value = 5528.643
$ws.Range("K136").Value = 11864.454
$ws.Range("L136").Value = 16585.929
$ws.Range("M136").Value = -9314.454000000002
$ws.Range("N136").Value = -21685.929

$ws.Range("H141").Value = 224423.19
$ws.Range("J141").Value = 236708.45
$ws.Range("L141").Value = 236708.45
$ws.Range("N141").Value = -247068.45


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13403640
$ws.Range("I4").Value = 9752329
$ws.Range("K4").Value = 29256987
$ws.Range("M4").Value = -29256875

$ws.Range("H10").Value = 59.5
$ws.Range("I10").Value = 73.333336
$ws.Range("K10").Value = 220.000008
$ws.Range("M10").Value = -81.00000800000001

$ws.Range("H82").Value = 5011666.5
$ws.Range("I82").Value = 5011666.5
$ws.Range("K82").Value = 15034999.5
$ws.Range("M82").Value = -15034593.5

$ws.Range("H85").Value = 5011666.5
$ws.Range("I85").Value = 5011666.5
$ws.Range("K85").Value = 15034999.5
$ws.Range("M85").Value = -15033595.5

$ws.Range("H121").Value = 5316621.5
$ws.Range("I121").Value = 831.5
$ws.Range("J121").Value = 7770063
$ws.Range("K121").Value = 2494.5
$ws.Range("L121").Value = 23310189
$ws.Range("M121").Value = -1184.5
$ws.Range("N121").Value = -23312809


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3207.1667
$ws.Range("I80").Value = 3098.25
$ws.Range("K80").Value = 3098.25
$ws.Range("M80").Value = -2100.25

$ws.Range("H83").Value = 3207.1667
$ws.Range("I83").Value = 3098.25
$ws.Range("K83").Value = 15491.25
$ws.Range("M83").Value = -10499.25

$ws.Range("H122").Value = 7411.1665
$ws.Range("I122").Value = 6117.125
$ws.Range("K122").Value = 18351.375
$ws.Range("M122").Value = -15901.375

$ws.Range("H132").Value = 4181.154
$ws.Range("J132").Value = 4402
$ws.Range("L132").Value = 13206
$ws.Range("N132").Value = -18266


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13788.179
$ws.Range("I7").Value = 10369.6875
$ws.Range("K7").Value = 10369.6875
$ws.Range("M7").Value = -10257.6875

$ws.Range("H22").Value = 10038.6
$ws.Range("I22").Value = 2600.5
$ws.Range("K22").Value = 2600.5
$ws.Range("M22").Value = -2305.5

$ws.Range("H27").Value = 10038.6
$ws.Range("I27").Value = 2600.5
$ws.Range("K27").Value = 2600.5
$ws.Range("M27").Value = -2493.5

$ws.Range("H46").Value = 3937.5
$ws.Range("I46").Value = 3978.261
$ws.Range("K46").Value = 3978.261
$ws.Range("M46").Value = -3790.261

$ws.Range("H126").Value = 13788.179
$ws.Range("I126").Value = 10369.6875
$ws.Range("K126").Value = 31109.0625
$ws.Range("M126").Value = -28639.0625

$ws.Range("H132").Value = 19530.967
$ws.Range("J132").Value = 19468.572
$ws.Range("L132").Value = 58405.716
$ws.Range("N132").Value = -63465.716

$ws.Range("H136").Value = 95662400
$ws.Range("I136").Value = 80009650
$ws.Range("J136").Value = 125011310
$ws.Range("K136").Value = 240028950
$ws.Range("L136").Value = 375033930
$ws.Range("M136").Value = -240026400
$ws.Range("N136").Value = -375039030


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4507.9
$ws.Range("I81").Value = 4507.9
$ws.Range("K81").Value = 9015.799999999999
$ws.Range("M81").Value = -7954.799999999999

$ws.Range("H84").Value = 4507.9
$ws.Range("I84").Value = 4507.9
$ws.Range("K84").Value = 45079
$ws.Range("M84").Value = -39775

$ws.Range("H107").Value = 33656.72
$ws.Range("I107").Value = 1881.2632
$ws.Range("K107").Value = 5643.7896
$ws.Range("M107").Value = -3723.7896

$ws.Range("H122").Value = 4132
$ws.Range("I122").Value = 2810.875
$ws.Range("K122").Value = 8432.625
$ws.Range("M122").Value = -5982.625

$ws.Range("H132").Value = 3529.7693
$ws.Range("I132").Value = 2394.675
$ws.Range("J132").Value = 7313.4165
$ws.Range("K132").Value = 7184.025000000001
$ws.Range("L132").Value = 21940.2495
$ws.Range("M132").Value = -4654.025000000001
$ws.Range("N132").Value = -27000.2495

$ws.Range("H136").Value = 4226.355
$ws.Range("I136").Value = 3005.9048
$ws.Range("K136").Value = 9017.714399999999
$ws.Range("M136").Value = -6467.714399999999

